# Update gh-pages to output generated at 456a3b4
# Refreshes the "想去人数" (interested-count, column F) figures scraped
# from bilibili show listings across the 展览 / 演出 / 全部类型 sheets.
# (本地生活 has no data rows, so it needs no changes.)

$wb = $excel.ActiveWorkbook

# --- 展览 (Exhibitions) ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 20412
$ws.Range("F6").Value = 1109
$ws.Range("F8").Value = 7662
$ws.Range("F9").Value = 526
$ws.Range("F10").Value = 745
$ws.Range("F11").Value = 283
$ws.Range("F12").Value = 46
$ws.Range("F14").Value = 133
$ws.Range("F15").Value = 21
$ws.Range("F17").Value = 202
$ws.Range("F18").Value = 1348
$ws.Range("F19").Value = 459
$ws.Range("F21").Value = 692
$ws.Range("F22").Value = 52
$ws.Range("F24").Value = 73
$ws.Range("F25").Value = 330
$ws.Range("F26").Value = 1129
$ws.Range("F28").Value = 23
$ws.Range("F30").Value = 5218
$ws.Range("F32").Value = 92
$ws.Range("F33").Value = 4892
$ws.Range("F34").Value = 27
$ws.Range("F35").Value = 93
$ws.Range("F37").Value = 12763
$ws.Range("F38").Value = 1342
$ws.Range("F40").Value = 35
$ws.Range("F43").Value = 389
$ws.Range("F44").Value = 4016
$ws.Range("F45").Value = 323

# --- 演出 (Performances) ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 215

# --- 全部类型 (All types, aggregated view) ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 20412
$ws.Range("F6").Value = 1109
$ws.Range("F8").Value = 7662
$ws.Range("F9").Value = 526
$ws.Range("F10").Value = 745
$ws.Range("F11").Value = 283
$ws.Range("F12").Value = 46
$ws.Range("F14").Value = 133
$ws.Range("F15").Value = 21
$ws.Range("F17").Value = 202
$ws.Range("F18").Value = 1348
$ws.Range("F19").Value = 459
$ws.Range("F21").Value = 692
$ws.Range("F22").Value = 52
$ws.Range("F24").Value = 73
$ws.Range("F25").Value = 330
$ws.Range("F26").Value = 1129
$ws.Range("F28").Value = 23
$ws.Range("F30").Value = 215
$ws.Range("F31").Value = 5218
$ws.Range("F34").Value = 92
$ws.Range("F36").Value = 4892
$ws.Range("F37").Value = 27
$ws.Range("F38").Value = 93
$ws.Range("F40").Value = 12763
$ws.Range("F41").Value = 1342
$ws.Range("F43").Value = 35
$ws.Range("F46").Value = 389
$ws.Range("F47").Value = 4016
$ws.Range("F48").Value = 323
